$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("C1").Value = "monosaccharides"
$ws.Range("D1").Value = "motifs"
$ws.Range("E1").Value = "sasa"
$ws.Range("G1").Value = "has_multi_node_motifs"

# --- Capture old J/K (monosaccharides/motifs) text values before overwriting C/D ---
$mono = @{}
$motifs = @{}
for ($r = 2; $r -le 6; $r++) {
    $mono[$r] = $ws.Cells.Item($r, 10).Value2  # column J
    $motifs[$r] = $ws.Cells.Item($r, 11).Value2 # column K
}

# --- New per-row values ---
$newFlex = @{
    2 = 0.234424633348246
    3 = 2.332801287369361
    4 = 0.7661902446741994
    5 = 0.8630683869676127
    6 = 0.2143102970853854
}

for ($r = 2; $r -le 6; $r++) {
    # C: monosaccharides (was column J)
    $ws.Cells.Item($r, 3).Value = $mono[$r]
    # D: motifs (was column K)
    $ws.Cells.Item($r, 4).Value = $motifs[$r]
    # E (sasa) keeps its existing max_SASA numeric value - no change needed
    # F: flexibility - replace with new aggregated value
    $ws.Cells.Item($r, 6).Value = $newFlex[$r]
    # G: has_multi_node_motifs boolean
    $ws.Cells.Item($r, 7).Value = $false
}

# --- Remove now-unused columns H through L ---
$ws.Range("H1:L6").Delete() | Out-Null

# --- Fix the sheet dimension ---
$ws.UsedRange | Out-Null
